$wb = $excel.ActiveWorkbook

# Data for the new row (r=42) on each of the 5 sheets.
# Column order: A=Date(serial), B..H = New York, Texas, California, Florida, Kansas, Nebraska, UK
$rowsData = @{
    1 = @(43944, 20973, 604, 1531, 987, 113, 38, 18738)
    2 = @(43944, 44381.54372776611, 1041.527243128084, 1937.375176284058, 2297.728107947313, 1939.372137709838, 982.212645936121, 13873.80077473827)
    3 = @(43944, 1560, 48, 112, 94, 1, 20, 638)
    4 = @(43944, 3301.159024236644, 82.77037693733121, 141.7282950645424, 218.8312483759346, 17.1625852894676, 516.9540241769058, 472.3815185336223)
    5 = @(43944, 1845.686474961026, 44.14420103324331, 98.95672030399302, 111.2780177911668, 96.11047762101857, 118.8994255606883, 484.8204049150721)
}

$cols = @("A", "B", "C", "D", "E", "F", "G", "H")

for ($i = 1; $i -le 5; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Duplicate the formatting of row 41 into row 42 (keeps the date-style
    # xf used on column A, e.g. s="2", instead of minting a new style).
    $ws.Range("A41:H41").Copy($ws.Range("A42:H42"))

    $values = $rowsData[$i]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + "42").Value = $values[$c]
    }
}
